# Add 2022 mortality data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several values in this table look numeric (e.g. "0.12", "104.00") but
# must be stored as text, matching the existing convention in this worksheet
# where every cell (including numeric-looking ones) is a shared string.
# Prefixing with a leading apostrophe forces Excel to keep the literal text;
# re-applying the "Normal" style afterwards clears the quote-prefix cell
# formatting so the cell itself carries no extra style reference.

# Hepatitis A: extend post-vaccine period to 2022 and update annual post-vaccine deaths
$ws.Range("D2").Value = "2007-2022"
$ws.Range("E2").Value = "'0.12"
$ws.Range("E2").Style = "Normal"

# Hepatitis B: extend post-vaccine period to 2022 and update annual post-vaccine deaths
$ws.Range("D3").Value = "2010-2022"
$ws.Range("E3").Value = "'0.15"
$ws.Range("E3").Style = "Normal"

# Meningococcus: update annual pre-vaccine deaths, extend post-vaccine period to 2022,
# and update annual post-vaccine deaths
$ws.Range("C4").Value = "'0.00"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "2018-2022"
$ws.Range("E4").Value = "'0.00"
$ws.Range("E4").Style = "Normal"

# Varicella: extend post-vaccine period to 2022 and update annual post-vaccine deaths
$ws.Range("D5").Value = "2016-2022"
$ws.Range("E5").Value = "'3.29"
$ws.Range("E5").Style = "Normal"

# Covid: extend period to 2022 and update annual deaths
$ws.Range("B6").Value = "2020-2022"
$ws.Range("C6").Value = "'104.00"
$ws.Range("C6").Style = "Normal"
